# Fruta / hortaliza, semanal
# Weekly refresh of the Maracuya (Vega Central Mapocho de Santiago) price series:
# existing daily records are updated with the newest week's figures and two
# older historical rows are appended at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh rows 4-30 with the updated weekly price data ---
$ws.Cells.Item(4, 4).Value = 44424
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 30
$ws.Cells.Item(4, 14).Value = 32000
$ws.Cells.Item(4, 15).Value = 32000
$ws.Cells.Item(4, 16).Value = 32000
$ws.Cells.Item(4, 19).Value = 1778

$ws.Cells.Item(5, 12).Value = "Especial"
$ws.Cells.Item(5, 13).Value = 25
$ws.Cells.Item(5, 14).Value = 33000
$ws.Cells.Item(5, 15).Value = 33000
$ws.Cells.Item(5, 16).Value = 33000
$ws.Cells.Item(5, 19).Value = 1833

$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 45
$ws.Cells.Item(6, 14).Value = 30000
$ws.Cells.Item(6, 15).Value = 30000
$ws.Cells.Item(6, 16).Value = 30000
$ws.Cells.Item(6, 19).Value = 1667

$ws.Cells.Item(7, 4).Value = 44403
$ws.Cells.Item(7, 12).Value = "Segunda"
$ws.Cells.Item(7, 13).Value = 15
$ws.Cells.Item(7, 14).Value = 28000
$ws.Cells.Item(7, 15).Value = 28000
$ws.Cells.Item(7, 16).Value = 28000
$ws.Cells.Item(7, 19).Value = 1556

$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 35
$ws.Cells.Item(8, 14).Value = 37000
$ws.Cells.Item(8, 15).Value = 37000
$ws.Cells.Item(8, 16).Value = 37000
$ws.Cells.Item(8, 19).Value = 2056

$ws.Cells.Item(9, 4).Value = 44396
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 15
$ws.Cells.Item(9, 14).Value = 34000
$ws.Cells.Item(9, 15).Value = 34000
$ws.Cells.Item(9, 16).Value = 34000
$ws.Cells.Item(9, 19).Value = 1889

$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 35
$ws.Cells.Item(10, 14).Value = 37000
$ws.Cells.Item(10, 15).Value = 37000
$ws.Cells.Item(10, 16).Value = 37000
$ws.Cells.Item(10, 19).Value = 2056

$ws.Cells.Item(11, 4).Value = 44340
$ws.Cells.Item(11, 12).Value = "Segunda"
$ws.Cells.Item(11, 13).Value = 20
$ws.Cells.Item(11, 14).Value = 35000
$ws.Cells.Item(11, 15).Value = 35000
$ws.Cells.Item(11, 16).Value = 35000
$ws.Cells.Item(11, 19).Value = 1944

$ws.Cells.Item(12, 4).Value = 44354
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 40
$ws.Cells.Item(12, 14).Value = 38000
$ws.Cells.Item(12, 15).Value = 38000
$ws.Cells.Item(12, 16).Value = 38000
$ws.Cells.Item(12, 19).Value = 2111

$ws.Cells.Item(13, 12).Value = "Especial"
$ws.Cells.Item(13, 13).Value = 15
$ws.Cells.Item(13, 14).Value = 32000
$ws.Cells.Item(13, 15).Value = 32000
$ws.Cells.Item(13, 16).Value = 32000
$ws.Cells.Item(13, 19).Value = 1778

$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 25
$ws.Cells.Item(14, 14).Value = 30000
$ws.Cells.Item(14, 15).Value = 30000
$ws.Cells.Item(14, 16).Value = 30000
$ws.Cells.Item(14, 19).Value = 1667

$ws.Cells.Item(15, 4).Value = 44410
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 10
$ws.Cells.Item(15, 14).Value = 28000
$ws.Cells.Item(15, 15).Value = 28000
$ws.Cells.Item(15, 16).Value = 28000
$ws.Cells.Item(15, 19).Value = 1556

$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 35
$ws.Cells.Item(16, 14).Value = 29000
$ws.Cells.Item(16, 15).Value = 29000
$ws.Cells.Item(16, 16).Value = 29000
$ws.Cells.Item(16, 19).Value = 1611

$ws.Cells.Item(17, 4).Value = 44389
$ws.Cells.Item(17, 12).Value = "Segunda"
$ws.Cells.Item(17, 13).Value = 20
$ws.Cells.Item(17, 14).Value = 27000
$ws.Cells.Item(17, 15).Value = 27000
$ws.Cells.Item(17, 16).Value = 27000
$ws.Cells.Item(17, 19).Value = 1500

$ws.Cells.Item(18, 4).Value = 44305
$ws.Cells.Item(18, 13).Value = 25
$ws.Cells.Item(18, 14).Value = 30000
$ws.Cells.Item(18, 15).Value = 30000
$ws.Cells.Item(18, 16).Value = 30000
$ws.Cells.Item(18, 19).Value = 1667

$ws.Cells.Item(19, 4).Value = 44417
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 15
$ws.Cells.Item(19, 14).Value = 28000
$ws.Cells.Item(19, 15).Value = 28000
$ws.Cells.Item(19, 16).Value = 28000
$ws.Cells.Item(19, 19).Value = 1556

$ws.Cells.Item(20, 4).Value = 44333
$ws.Cells.Item(20, 13).Value = 30
$ws.Cells.Item(20, 14).Value = 38000
$ws.Cells.Item(20, 15).Value = 38000
$ws.Cells.Item(20, 16).Value = 38000
$ws.Cells.Item(20, 19).Value = 2111

$ws.Cells.Item(21, 4).Value = 44333
$ws.Cells.Item(21, 13).Value = 25
$ws.Cells.Item(21, 14).Value = 35000
$ws.Cells.Item(21, 15).Value = 35000
$ws.Cells.Item(21, 16).Value = 35000
$ws.Cells.Item(21, 19).Value = 1944

$ws.Cells.Item(22, 4).Value = 44277
$ws.Cells.Item(22, 13).Value = 100
$ws.Cells.Item(22, 14).Value = 30000
$ws.Cells.Item(22, 15).Value = 30000
$ws.Cells.Item(22, 16).Value = 30000
$ws.Cells.Item(22, 19).Value = 1667

$ws.Cells.Item(23, 4).Value = 44277
$ws.Cells.Item(23, 12).Value = "Segunda"
$ws.Cells.Item(23, 13).Value = 60
$ws.Cells.Item(23, 14).Value = 28000
$ws.Cells.Item(23, 15).Value = 28000
$ws.Cells.Item(23, 16).Value = 28000
$ws.Cells.Item(23, 19).Value = 1556

$ws.Cells.Item(24, 4).Value = 44319
$ws.Cells.Item(24, 13).Value = 140
$ws.Cells.Item(24, 14).Value = 27000
$ws.Cells.Item(24, 15).Value = 27000
$ws.Cells.Item(24, 16).Value = 27000
$ws.Cells.Item(24, 19).Value = 1500

$ws.Cells.Item(25, 12).Value = "Especial"
$ws.Cells.Item(25, 13).Value = 16
$ws.Cells.Item(25, 14).Value = 35000
$ws.Cells.Item(25, 15).Value = 35000
$ws.Cells.Item(25, 16).Value = 35000
$ws.Cells.Item(25, 19).Value = 1944

$ws.Cells.Item(26, 4).Value = 44326
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 25
$ws.Cells.Item(26, 14).Value = 30000
$ws.Cells.Item(26, 15).Value = 30000
$ws.Cells.Item(26, 16).Value = 30000
$ws.Cells.Item(26, 19).Value = 1667

$ws.Cells.Item(27, 4).Value = 44326
$ws.Cells.Item(27, 12).Value = "Segunda"
$ws.Cells.Item(27, 13).Value = 20
$ws.Cells.Item(27, 14).Value = 28000
$ws.Cells.Item(27, 15).Value = 28000
$ws.Cells.Item(27, 16).Value = 28000
$ws.Cells.Item(27, 19).Value = 1556

$ws.Cells.Item(28, 12).Value = "Especial"
$ws.Cells.Item(28, 13).Value = 20
$ws.Cells.Item(28, 14).Value = 35000
$ws.Cells.Item(28, 15).Value = 35000
$ws.Cells.Item(28, 16).Value = 35000
$ws.Cells.Item(28, 19).Value = 1944

$ws.Cells.Item(29, 4).Value = 44382
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 30
$ws.Cells.Item(29, 14).Value = 32000
$ws.Cells.Item(29, 15).Value = 32000
$ws.Cells.Item(29, 16).Value = 32000
$ws.Cells.Item(29, 19).Value = 1778

$ws.Cells.Item(30, 4).Value = 44382
$ws.Cells.Item(30, 12).Value = "Segunda"
$ws.Cells.Item(30, 13).Value = 15
$ws.Cells.Item(30, 14).Value = 30000
$ws.Cells.Item(30, 15).Value = 30000
$ws.Cells.Item(30, 16).Value = 30000
$ws.Cells.Item(30, 19).Value = 1667

# --- Append two additional historical rows (31-32) ---
$ws.Cells.Item(31, 1).Value = 9
$ws.Cells.Item(31, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(31, 3).Value = "Metropolitana"
$ws.Cells.Item(31, 4).Value = 44270
$ws.Cells.Item(31, 5).Value = 13
$ws.Cells.Item(31, 6).Value = "Fruta"
$ws.Cells.Item(31, 7).Value = 100108
$ws.Cells.Item(31, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(31, 9).Value = 100108003
$ws.Cells.Item(31, 10).Value = "Maracuyá"
$ws.Cells.Item(31, 11).Value = "Sin especificar"
$ws.Cells.Item(31, 12).Value = "Especial"
$ws.Cells.Item(31, 13).Value = 70
$ws.Cells.Item(31, 14).Value = 38000
$ws.Cells.Item(31, 15).Value = 38000
$ws.Cells.Item(31, 16).Value = 38000
$ws.Cells.Item(31, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(31, 18).Value = "Perú"
$ws.Cells.Item(31, 19).Value = 2111
$ws.Cells.Item(31, 20).Value = 18
$ws.Cells.Item(31, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

$ws.Cells.Item(32, 1).Value = 9
$ws.Cells.Item(32, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(32, 3).Value = "Metropolitana"
$ws.Cells.Item(32, 4).Value = 44284
$ws.Cells.Item(32, 5).Value = 13
$ws.Cells.Item(32, 6).Value = "Fruta"
$ws.Cells.Item(32, 7).Value = 100108
$ws.Cells.Item(32, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(32, 9).Value = 100108003
$ws.Cells.Item(32, 10).Value = "Maracuyá"
$ws.Cells.Item(32, 11).Value = "Sin especificar"
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 120
$ws.Cells.Item(32, 14).Value = 23000
$ws.Cells.Item(32, 15).Value = 23000
$ws.Cells.Item(32, 16).Value = 23000
$ws.Cells.Item(32, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(32, 18).Value = "Perú"
$ws.Cells.Item(32, 19).Value = 1278
$ws.Cells.Item(32, 20).Value = 18
$ws.Cells.Item(32, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
